$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RES installed")

# Update the installed RES capacities (Pinst, [MW]) for nodes 4,5,6,7,8.
# This cascades through all dependent formulas across the workbook
# (Main!B7 total, and all the Pg/Pc profile sheets that VLOOKUP into
# this table) on recalculation.
$ws.Range("C2").Value = 50
$ws.Range("C3").Value = 50
$ws.Range("C4").Value = 40
$ws.Range("C5").Value = 40
$ws.Range("C6").Value = 40

# Make "RES installed" the active sheet/tab, with cell R13 selected,
# matching the author's final view state when the workbook was saved.
$ws.Activate()
$ws.Range("R13").Select()
